$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4 (everything from row 4 down shifts by one)
$ws.Rows("4").Insert()

# Move the small "legend" cells that used to sit at H1 / I1 down into
# column A (rows 4 and 5), carrying over their fill colours.
$ws.Range("A4").Value = $ws.Range("H1").Value2
$ws.Range("A4").Interior.Color = $ws.Range("H1").Interior.Color

$ws.Range("A5").Value = $ws.Range("I1").Value2
$ws.Range("A5").Interior.Color = $ws.Range("I1").Interior.Color

# Remove the now-empty legend cells from row 1 entirely.
$ws.Range("H1").Clear()
$ws.Range("I1").Clear()

# Widen column A.
$ws.Columns("A").ColumnWidth = 31.16666666666667

# Update status text that changed in the dataset rows.
$ws.Range("H9").Value = "ready to be fit"
$ws.Range("I9").Value = "ready to be fit"
$ws.Range("K9").Value = "ready to be fit"

$ws.Range("H10").Value = "ready to be fit"
$ws.Range("J10").Value = "ready to be fit"

$ws.Range("C11").Value = "hadding"

# Restore the previously-selected cell (shifted down with the insert).
$ws.Range("D20").Select()
